$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing G-column (Mean) values for rows 4,5,6,7,16,17,18,19,32,33,34,35
$gUpdates = @{
    4 = 0.253762984106252
    5 = 0.253762984106252
    6 = 0.573227640253274
    7 = 0.573227640253274
    16 = 0.271142099673189
    17 = 0.271142099673189
    18 = 0.578541131053511
    19 = 0.578541131053511
    32 = 0.239414178053617
    33 = 0.239414178053617
    34 = 0.601740116071675
    35 = 0.601740116071675
}

foreach ($r in $gUpdates.Keys) {
    $ws.Cells.Item($r, 7).Value = $gUpdates[$r]
}

# New rows 42-57 (2019 - 2023 dataset block)
$rows = @(
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'DRP (95th Percentile)'; 'C' = 'D'; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 0.0785; 'G' = 0.111948275862069; 'H' = 0.946; 'I' = 0.2128; 'J' = $null; 'K' = $null; 'L' = 0.1555; 'M' = 0.16448; 'N' = 0.2046; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'mg/L' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'DRP (Median)'; 'C' = 'D'; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 0.0785; 'G' = 0.111948275862069; 'H' = 0.946; 'I' = 0.2128; 'J' = $null; 'K' = $null; 'L' = 0.1555; 'M' = 0.16448; 'N' = 0.2046; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'mg/L' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'E coli (>260)'; 'C' = 'D'; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 262.5; 'G' = 518.163838671249; 'H' = 9700; 'I' = 1120; 'J' = 18.9655172413793; 'K' = 50; 'L' = 201.5; 'M' = 570; 'N' = 886; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = '% exceedances over 260/100 mL' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'E coli (>540)'; 'C' = 'C'; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 262.5; 'G' = 518.163838671249; 'H' = 9700; 'I' = 1120; 'J' = 18.9655172413793; 'K' = 50; 'L' = 201.5; 'M' = 570; 'N' = 886; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = '% exceedances over 540/100 mL' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'E coli (Median)'; 'C' = 'E'; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 262.5; 'G' = 518.163838671249; 'H' = 9700; 'I' = 1120; 'J' = 18.9655172413793; 'K' = 50; 'L' = 201.5; 'M' = 570; 'N' = 886; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'E. coli/100 mL' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'E coli (95th Percentile)'; 'C' = 'C'; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 262.5; 'G' = 518.163838671249; 'H' = 9700; 'I' = 1120; 'J' = 18.9655172413793; 'K' = 50; 'L' = 201.5; 'M' = 570; 'N' = 886; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'E. coli/100 mL' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'Ammoniacal-N (95th Percentile)'; 'C' = 'C'; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 0.12609; 'G' = 0.184900156157885; 'H' = 0.83691117624023; 'I' = 0.60837; 'J' = $null; 'K' = $null; 'L' = 0.0783; 'M' = 0.34424; 'N' = 0.44129; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'mg NH4-N/L' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'Ammoniacal-N (Median)'; 'C' = 'B'; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 0.12609; 'G' = 0.184900156157885; 'H' = 0.83691117624023; 'I' = 0.60837; 'J' = $null; 'K' = $null; 'L' = 0.0783; 'M' = 0.34424; 'N' = 0.44129; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'mg NH4-N/L' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'Nitrate-N (95th Percentile)'; 'C' = 'A'; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 0.608; 'G' = 0.57744650083364; 'H' = 1.76; 'I' = 1.086; 'J' = $null; 'K' = $null; 'L' = 0.305; 'M' = 0.88296; 'N' = 0.99288; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'mg NO3-N/L' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'Nitrate-N (Median)'; 'C' = 'A'; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 0.608; 'G' = 0.57744650083364; 'H' = 1.76; 'I' = 1.086; 'J' = $null; 'K' = $null; 'L' = 0.305; 'M' = 0.88296; 'N' = 0.99288; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'mg NO3-N/L' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'Soluble Inorganic Nitrogen (95th Percentile)'; 'C' = $null; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 0.88; 'G' = 0.881724137931034; 'H' = 2.26; 'I' = 1.748; 'J' = $null; 'K' = $null; 'L' = 0.52; 'M' = 1.3264; 'N' = 1.4458; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'g/m3' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'Soluble Inorganic Nitrogen (Median)'; 'C' = $null; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 0.88; 'G' = 0.881724137931034; 'H' = 2.26; 'I' = 1.748; 'J' = $null; 'K' = $null; 'L' = 0.52; 'M' = 1.3264; 'N' = 1.4458; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'g/m3' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'Total Nitrogen (95th Percentile)'; 'C' = $null; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 2.21; 'G' = 2.26396551724138; 'H' = 3.57; 'I' = 3.234; 'J' = $null; 'K' = $null; 'L' = 2.12; 'M' = 2.81; 'N' = 3.0944; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'g/m3' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'Total Nitrogen (Median)'; 'C' = $null; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 2.21; 'G' = 2.26396551724138; 'H' = 3.57; 'I' = 3.234; 'J' = $null; 'K' = $null; 'L' = 2.12; 'M' = 2.81; 'N' = 3.0944; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'g/m3' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'Total Phosphorus (95th Percentile)'; 'C' = $null; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 0.2525; 'G' = 0.285068965517241; 'H' = 1.23; 'I' = 0.5454; 'J' = $null; 'K' = $null; 'L' = 0.37; 'M' = 0.37028; 'N' = 0.50094; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'g/m3' }
    @{ 'A' = 'Whitebait Creek at Edinburgh Terrace'; 'B' = 'Total Phosphorus (Median)'; 'C' = $null; 'D' = '2019 - 2023'; 'E' = 'RepSite'; 'F' = 0.2525; 'G' = 0.285068965517241; 'H' = 1.23; 'I' = 0.5454; 'J' = $null; 'K' = $null; 'L' = 0.37; 'M' = 0.37028; 'N' = 0.50094; 'O' = 1790056; 'P' = 5517935; 'Q' = 'Horowhenua District'; 'R' = 'Manawatū'; 'S' = 'Coastal Manawatu'; 'T' = 'Mana_13a'; 'U' = 'g/m3' }
)

$startRow = 42
$colIndex = @{ 'A'=1; 'B'=2; 'C'=3; 'D'=4; 'E'=5; 'F'=6; 'G'=7; 'H'=8; 'I'=9; 'J'=10; 'K'=11; 'L'=12; 'M'=13; 'N'=14; 'O'=15; 'P'=16; 'Q'=17; 'R'=18; 'S'=19; 'T'=20; 'U'=21 }

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    foreach ($col in $rowData.Keys) {
        $c = $colIndex[$col]
        $val = $rowData[$col]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}


